$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.77%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.68%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.155"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.28%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08094"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'9.36%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.475"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'11.72%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.796"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.76%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.915"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.39%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9288"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.06%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1762"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.97%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07424"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.82%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08848"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'8.90%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03010"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1001"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.84%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001529"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'2.44%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006015"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.08%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.528"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.59%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'2.97%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.31%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1339"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.48%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.156"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-10.62%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1681"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'7.24%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04622"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.31%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001240"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.14%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004529"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.24%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-7.59%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003411"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.51%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01749"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.78%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04609"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.40%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006925"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-5.57%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1372"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.88%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002189"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.68%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-2.97%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006206"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.37%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.008402"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-15.91%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.7485"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-8.85%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("E50").Style = "Normal"
